# Fruta / hortaliza, semanal
# Update the weekly price records for "Zapallo italiano" (Mapocho Venta Directa
# de Santiago) by re-aligning each row's data with the correct week it belongs
# to (the rows had been shifted out of their proper weekly order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = 44243
    "J2"  = 80
    "L2"  = 11000
    "M2"  = 10375
    "O2"  = "Provincia de Quillota"
    "P2"  = 173

    "D3"  = 44585
    "J3"  = 30
    "K3"  = 11000
    "M3"  = 11000
    "O3"  = "Provincia de Limarí"
    "P3"  = 183

    "D4"  = 44315
    "J4"  = 25
    "K4"  = 10000
    "L4"  = 10000
    "M4"  = 10000
    "N4"  = "`$/caja 60 unidades"
    "O4"  = "Provincia de Limarí"
    "P4"  = 167
    "Q4"  = 60

    "D5"  = 44277
    "J5"  = 25
    "K5"  = 10000
    "L5"  = 10000
    "M5"  = 10000
    "P5"  = 167

    "D6"  = 44291
    "J6"  = 20

    "D7"  = 44186
    "J7"  = 15
    "K7"  = 7000
    "L7"  = 7000
    "M7"  = 7000
    "P7"  = 117

    "D9"  = 44200
    "J9"  = 10

    "D10" = 44405
    "J10" = 45
    "K10" = 9000
    "L10" = 9000
    "M10" = 9000
    "N10" = "`$/caja 50 unidades"
    "O10" = "Provincia de Quillota"
    "P10" = 180
    "Q10" = 50

    "D11" = 44312
    "J11" = 30
    "K11" = 10000
    "L11" = 10000
    "M11" = 10000
    "P11" = 167

    "D12" = 44284
    "J12" = 35
    "K12" = 10000
    "L12" = 10000
    "M12" = 10000
    "P12" = 167

    "D13" = 44179
    "J13" = 15
    "K13" = 7000
    "L13" = 7000
    "M13" = 7000
    "P13" = 117
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
